# "Improved help and tests"
#
# The three instruction rows (2-4) used to hold a single long string in
# column A ("Label:  full command line ..."). The edit splits each row into
# a short label in column A and the actual command line in column B, and
# tidies up the sample command lines themselves (no more "=" in the long
# options, "--merge Color" / "--append" spelled out properly).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Diff test:"
$ws.Range("B2").Value = "xltablediff.py  --key ID test1old.xlsx test1new.xlsx --out test1diff.xlsx"

$ws.Range("A3").Value = "Merge test:"
$ws.Range("B3").Value = "xltablediff.py  --key ID --merge Color test1old.xlsx test1new.xlsx --out test1merge.xlsx"

$ws.Range("A4").Value = "Append test:"
$ws.Range("B4").Value = "xltablediff.py  --key ID --append test1old.xlsx test1new.xlsx --out test1append.xlsx"

# The new command-line cells in column B pick up a plain Arial 10 (regular)
# font, distinct from the label cells.
$ws.Range("B2:B4").Font.Name = "Arial"
$ws.Range("B2:B4").Font.Size = 10
$ws.Range("B2:B4").Font.Bold = $false
$ws.Range("B2:B4").Font.Italic = $false

# Selection ends up covering the newly split label/command cells.
$ws.Range("A2:B4").Select()
